# Daily attendance processing - 2026-01-25 17:34:14
#
# The "Recorded By" column (G) lists the users who recorded/edited each
# attendance session, as a comma-separated string. This pass rotates that
# list so the most recent recorder (previously trailing) moves to the
# front - except rows whose list includes "admin@admin.com", which are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($null -eq $raw) {
        continue
    }

    $text = [string]$raw

    if ($text -notlike "*,*") {
        continue
    }

    if ($text -like "*admin@admin.com*") {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Count -lt 2) {
        continue
    }

    $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
    $newText = [string]::Join(", ", $rotated)

    $cell.Value2 = $newText
}
